{"js": "// \"alterar cor no doc\" \u2014 three paragraphs in the body were left in the\n// default/black color while every sibling paragraph around them already\n// uses red (FF0000) Arial 12pt text. Bring those three into line by\n// applying the same red font color to them (paragraph mark included, so\n// the <w:pPr>/<w:rPr> and every run's <w:rPr> both get <w:color val=\"FF0000\"/>).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = new Set([\n  \"Salvar como\u2026 (op\u00e7\u00e3o de escolher o nome do arquivo)\",\n  \"Corre\u00e7\u00e3o gamma... (valor)\",\n  \"Separa\u00e7\u00e3o de camadas R, G e B. (uma de cada vez, uma op\u00e7\u00e3o para cada)\"\n]);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text.trim();\n  if (targets.has(text)) {\n    // Setting color on the paragraph's font applies it to the whole\n    // paragraph range, including the paragraph mark, matching the diff.\n    paragraph.font.color = \"#FF0000\";\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraphs whose text (trimmed) should become red (RGB FF0000), matching\n# the \"alterar cor no doc\" commit: these three paragraphs were still plain\n# black/default-colored while their siblings already used red Arial 12pt.\n$targets = @(\n    \"Salvar como\u2026 (op\u00e7\u00e3o de escolher o nome do arquivo)\",\n    \"Corre\u00e7\u00e3o gamma... (valor)\",\n    \"Separa\u00e7\u00e3o de camadas R, G e B. (uma de cada vez, uma op\u00e7\u00e3o para cada)\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($targets -contains $text) {\n        # 255 == wdColorRed, stored as RRGGBB \"FF0000\" in the OOXML <w:color>.\n        # Applying it to the paragraph's Range (not just the run) also colors\n        # the paragraph mark, so <w:pPr>/<w:rPr> picks up <w:color> too.\n        $p.Range.Font.Color = 255\n    }\n}\n"}
